# Adds a "Nome_estado" column (full state name) between the UF-abbreviation
# column (B) and the REGIAO column (old C, now D). Also renames the header
# of column B from "UF" to "Estado".
#
# Original layout: A=CODUF B=UF     C=REGIAO      D=CODREG
# New layout:       A=CODUF B=Estado C=Nome_estado D=REGIAO E=CODREG

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before current column C (REGIAO) so it becomes column D.
$ws.Columns.Item(3).Insert()

# Header for the state-abbreviation column was renamed from "UF" to "Estado".
$ws.Range("B1").Value = "Estado"

$ws.Range("C1").Value = "Nome_estado"

$names = @(
    "Rondônia",
    "Acre",
    "Amazonas",
    "Roraima",
    "Pará",
    "Amapá",
    "Tocantins",
    "Maranhão",
    "Piauí",
    "Ceará",
    "Rio Grande do Norte",
    "Paraíba",
    "Pernambuco",
    "Alagoas",
    "Sergipe",
    "Bahia",
    "Minas gerais",
    "Espiríto Santo",
    "Rio de Janeiro",
    "São Paulo",
    "Paraná",
    "Santa Catarina",
    "Rio Grande do Sul",
    "Mato Grosso do Sul",
    "Mato Grosso",
    "Goiás",
    "Distrito Federal"
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $names[$i]
}

$ws.Columns.Item(3).EntireColumn.AutoFit()

# Mirror the author's cursor position after typing the last name and
# pressing Enter (one row below the last filled cell, same column).
$ws.Range("C29").Select() | Out-Null
